$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3150
$ws.Range("I86").Value = 2866.6667
$ws.Range("K86").Value = 2866.6667
$ws.Range("M86").Value = -1743.6667
$ws.Range("H89").Value = 3150
$ws.Range("I89").Value = 2866.6667
$ws.Range("K89").Value = 14333.3335
$ws.Range("M89").Value = -8717.3335
$ws.Range("H116").Value = 20003450
$ws.Range("I116").Value = 33337350
$ws.Range("J116").Value = 2598.5
$ws.Range("K116").Value = 33337350
$ws.Range("L116").Value = 2598.5
$ws.Range("M116").Value = -33333908
$ws.Range("N116").Value = -9482.5
$ws.Range("H129").Value = 1361.7333
$ws.Range("J129").Value = 1776.2222
$ws.Range("L129").Value = 5328.6666
$ws.Range("N129").Value = -15328.6666
$ws.Range("H138").Value = 3985.5
$ws.Range("I138").Value = 3096.1277
$ws.Range("J138").Value = 4774.1885
$ws.Range("K138").Value = 9288.3831
$ws.Range("L138").Value = 14322.5655
$ws.Range("M138").Value = -4148.383099999999
$ws.Range("N138").Value = -24602.5655
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 45000
$ws.Range("J24").Value = 45000
$ws.Range("L24").Value = 45000
$ws.Range("N24").Value = -45748
$ws.Range("H32").Value = 57544
$ws.Range("I32").Value = 63163.332
$ws.Range("J32").Value = 50800.8
$ws.Range("K32").Value = 63163.332
$ws.Range("L32").Value = 50800.8
$ws.Range("M32").Value = -62876.332
$ws.Range("N32").Value = -51374.8
$ws.Range("H45").Value = 17316.5
$ws.Range("I45").Value = 779.8
$ws.Range("J45").Value = 100000
$ws.Range("K45").Value = 779.8
$ws.Range("L45").Value = 100000
$ws.Range("M45").Value = -402.8
$ws.Range("N45").Value = -100754
$ws.Range("H100").Value = 45000
$ws.Range("J100").Value = 45000
$ws.Range("L100").Value = 45000
$ws.Range("N100").Value = -47164
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 823.8333
$ws.Range("I29").Value = 823.8333
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 823.8333
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -534.8333
$ws.Range("N29").ClearContents()
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H109").Value = 25832.666
$ws.Range("J109").Value = 25832.666
$ws.Range("L109").Value = 25832.666
$ws.Range("N109").Value = -28606.666
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3721.2415
$ws.Range("I22").Value = 4566.2607
$ws.Range("J22").Value = 482
$ws.Range("K22").Value = 4566.2607
$ws.Range("L22").Value = 482
$ws.Range("M22").Value = -4216.2607
$ws.Range("N22").Value = -1182
$ws.Range("H31").Value = 4691.6665
$ws.Range("I31").Value = 3374.8333
$ws.Range("K31").Value = 3374.8333
$ws.Range("M31").Value = -3079.8333
$ws.Range("H34").Value = 4691.6665
$ws.Range("I34").Value = 3374.8333
$ws.Range("K34").Value = 3374.8333
$ws.Range("M34").Value = -3172.8333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1845.5938
$ws.Range("J5").Value = 2204.8667
$ws.Range("L5").Value = 6614.6001
$ws.Range("N5").Value = -6838.6001
$ws.Range("H113").Value = 3216.422
$ws.Range("I113").Value = 670
$ws.Range("J113").Value = 3767
$ws.Range("K113").Value = 2010
$ws.Range("L113").Value = 11301
$ws.Range("M113").Value = 160
$ws.Range("N113").Value = -15641
$ws.Range("H129").Value = 2502251.8
$ws.Range("J129").Value = 2780174
$ws.Range("L129").Value = 8340522
$ws.Range("N129").Value = -8350522
$ws.Range("H131").Value = 18871256
$ws.Range("I131").Value = 20550
$ws.Range("J131").Value = 20834872
$ws.Range("K131").Value = 61650
$ws.Range("L131").Value = 62504616
$ws.Range("M131").Value = -56610
$ws.Range("N131").Value = -62514696
$ws.Range("H135").Value = 1845.5938
$ws.Range("J135").Value = 2204.8667
$ws.Range("L135").Value = 19843.8003
$ws.Range("N135").Value = -24913.8003
$ws.Range("H138").Value = 1574.375
$ws.Range("I138").Value = 1094.6154
$ws.Range("J138").Value = 3653.3333
$ws.Range("K138").Value = 3283.8462
$ws.Range("L138").Value = 10959.9999
$ws.Range("M138").Value = 1856.1538
$ws.Range("N138").Value = -21239.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 33303.535
$ws.Range("J123").Value = 33303.535
$ws.Range("L123").Value = 33303.535
$ws.Range("N123").Value = -38203.535
$ws.Range("H126").Value = 3665.261
$ws.Range("I126").Value = 2839.3
$ws.Range("J126").Value = 4300.615
$ws.Range("K126").Value = 8517.900000000001
$ws.Range("L126").Value = 12901.845
$ws.Range("M126").Value = -6047.900000000001
$ws.Range("N126").Value = -17841.845
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4683.8335
$ws.Range("I7").Value = 4700.4443
$ws.Range("J7").Value = 4667.222
$ws.Range("K7").Value = 4700.4443
$ws.Range("L7").Value = 4667.222
$ws.Range("M7").Value = -4588.4443
$ws.Range("N7").Value = -4891.222
$ws.Range("H97").Value = 24303.908
$ws.Range("J97").Value = 24303.908
$ws.Range("L97").Value = 24303.908
$ws.Range("N97").Value = -26285.908
$ws.Range("H126").Value = 4683.8335
$ws.Range("I126").Value = 4700.4443
$ws.Range("J126").Value = 4667.222
$ws.Range("K126").Value = 14101.3329
$ws.Range("L126").Value = 14001.666
$ws.Range("M126").Value = -11631.3329
$ws.Range("N126").Value = -18941.666
$ws.Range("H131").Value = 25326
$ws.Range("J131").Value = 25326
$ws.Range("L131").Value = 25326
$ws.Range("N131").Value = -35406
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 10000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10336
$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 50000
$ws.Range("K20").Value = 50000
$ws.Range("M20").Value = -49760
$ws.Range("H28").Value = 16162.5
$ws.Range("I28").Value = 26350
$ws.Range("J28").Value = 5975
$ws.Range("K28").Value = 26350
$ws.Range("L28").Value = 5975
$ws.Range("M28").Value = -26002
$ws.Range("N28").Value = -6671
$ws.Range("H123").Value = 24347.125
$ws.Range("J123").Value = 24347.125
$ws.Range("L123").Value = 24347.125
$ws.Range("N123").Value = -34147.125
